# "further cleaning to metadata"
#
# For every data row (2-41) on the active sheet:
#   - column K (libraryProtocol) gets one single, shared value "E7420"
#     instead of the old, mostly-unique "E776x" library-protocol ids.
#     Re-using one value also collapses the now-unused shared-string
#     entries when the workbook is saved. The column is also re-styled
#     onto a dedicated Arial 11 font.
#   - column L (roboticLibraryPrep) becomes a real formula, =FALSE(),
#     instead of a bare boolean literal.
# Finally, the sheet's active selection is moved from L2:L41 to K2:K41
# to match where the edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 41; $r++) {
    $kCell = $ws.Cells.Item($r, 11)
    $kCell.Value = "E7420"
    $kCell.Font.Name = "Arial"
    $kCell.Font.Size = 11
    $kCell.Font.Color = 0

    $lCell = $ws.Cells.Item($r, 12)
    $lCell.Formula = "=FALSE()"
}

$ws.Range("K2:K41").Select() | Out-Null
